$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coin "Price" column (D) holds free-form numeric-looking strings (e.g. "42.383.30",
# "0.630") that must stay literal text -- force text format before writing so Excel
# does not reinterpret/round them as numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "42.383.30"
$ws.Range("E2").Value = "  -0.99%  "
$ws.Range("D3").Value = "2.256.45"
$ws.Range("E3").Value = "  -1.18%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "248.13"
$ws.Range("E5").Value = "  -1.35%  "
$ws.Range("D6").Value = "0.630"
$ws.Range("E6").Value = "  -2.31%  "
$ws.Range("D7").Value = "74.46"
$ws.Range("E7").Value = "  -0.31%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").Value = "0.622"
$ws.Range("E9").Value = "  -2.56%  "
$ws.Range("D10").Value = "41.65"
$ws.Range("E10").Value = "  +4.49%  "
$ws.Range("D11").Value = "0.0945"
$ws.Range("E11").Value = "  -3.11%  "
$ws.Range("D12").Value = "7.08"
$ws.Range("E12").Value = "  -5.50%  "
$ws.Range("D14").Value = "2.590.15"
$ws.Range("E14").Value = "  -1.18%  "
$ws.Range("D15").Value = "14.57"
$ws.Range("E15").Value = "  -3.03%  "
$ws.Range("D16").Value = "0.856"
$ws.Range("E16").Value = "  -2.01%  "
$ws.Range("D17").Value = "2.249.13"
$ws.Range("E17").Value = "  -1.39%  "
$ws.Range("D18").Value = "42.225.19"
$ws.Range("E18").Value = "  -1.03%  "
$ws.Range("D19").Value = "0.0₃0981"
$ws.Range("E19").Value = "  -2.25%  "
$ws.Range("E20").Value = "  -1.51%  "
$ws.Range("D21").Value = "71.91"
$ws.Range("E21").Value = "  -0.81%  "
$ws.Range("E22").Value = "  +4.82%  "
$ws.Range("D23").Value = "231.01"
$ws.Range("E23").Value = "  -1.89%  "
$ws.Range("E24").Value = "  +0.07%  "
$ws.Range("D25").Value = "11.11"
$ws.Range("E25").Value = "  -2.06%  "
$ws.Range("D26").Value = "7.99"
$ws.Range("E26").Value = "  +26.46%  "
$ws.Range("D27").Value = "3.56"
$ws.Range("E27").Value = "  -8.01%  "
$ws.Range("D28").Value = "2.31"
$ws.Range("E28").Value = "  -3.40%  "
$ws.Range("E29").Value = "  +2.90%  "
$ws.Range("D30").Value = "169.30"
$ws.Range("E30").Value = "  +0.86%  "
$ws.Range("D31").Value = "20.22"
$ws.Range("E31").Value = "  -4.02%  "
$ws.Range("D32").Value = "0.0827"
$ws.Range("E32").Value = "  -6.34%  "
$ws.Range("D33").Value = "0.120"
$ws.Range("E33").Value = "  -5.65%  "
$ws.Range("D34").Value = "30.51"
$ws.Range("E34").Value = "  -3.80%  "
$ws.Range("E35").Value = "  -2.34%  "
$ws.Range("D36").Value = "4.51"
$ws.Range("E36").Value = "  -1.66%  "
$ws.Range("D37").Value = "4.90"
$ws.Range("E37").Value = "  +2.45%  "
$ws.Range("D38").Value = "0.0304"
$ws.Range("E38").Value = "  -0.75%  "
$ws.Range("D39").Value = "13.68"
$ws.Range("E39").Value = "  -0.47%  "
$ws.Range("D40").Value = "2.19"
$ws.Range("E40").Value = "  -5.04%  "
$ws.Range("D41").Value = "5.80"
$ws.Range("E41").Value = "  -1.61%  "
$ws.Range("D42").Value = "61.90"
$ws.Range("E42").Value = "  +0.92%  "
$ws.Range("D43").Value = "0.204"
$ws.Range("E43").Value = "  -3.15%  "
$ws.Range("D44").Value = "108.55"
$ws.Range("E44").Value = "  +3.39%  "
$ws.Range("D45").Value = "8.68"
$ws.Range("E45").Value = "  -3.67%  "
$ws.Range("E46").Value = "  -0.25%  "
$ws.Range("D47").Value = "0.998"
$ws.Range("E47").Value = "  -0.26%  "
$ws.Range("D48").Value = "1.13"
$ws.Range("E48").Value = "  -3.27%  "
$ws.Range("D49").Value = "1.17"
$ws.Range("E49").Value = "  -0.90%  "
$ws.Range("D50").Value = "2.29"
$ws.Range("E50").Value = "  +0.78%  "
$ws.Range("B51").Value = "SynthetixNetwork"
$ws.Range("C51").Value = "https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx"
$ws.Range("D51").Value = "4.12"
$ws.Range("E51").Value = "  -2.30%  "
